# Prefix each step/row "Name" value (column A) in the listed sheets with the
# sheet's own name, e.g. "Step4 Seed" -> "free1 Step4 Seed" on sheet "free1".
# Row 1 (the "Name" header) is left untouched; only data rows (2..N) change.

$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "free1", "free2",
    "nomoney1", "nomoney2",
    "noppv1", "noppv2",
    "card1", "card2",
    "nosex1", "nosex2",
    "offtopic1", "offtopic2",
    "real1", "real2",
    "voice1", "voice2",
    "customyes1", "customyes2",
    "customno1", "customno2",
    "done1", "done2",
    "cumcontrol",
    "dickpic",
    "boosters",
    "price1", "price2",
    "discount1", "discount2"
)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $current = $cell.Value()
        if ($current -ne $null -and $current -ne "") {
            $prefixed = $sheetName + " " + $current
            $cell.Value = $prefixed
        }
    }
}
